$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the mailto hyperlinks: drop D3's link, refresh D2 / D4 ---
# (Hyperlinks.Add also stamps the cell text with the display string, so
#  the full e-mail values are (re)written afterwards.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:murbina@yopmail.com", "", "", "murbina@yopmail")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jmagallanes@yopmail.com", "", "", "jmagallanes@yopmail")

# --- Update the user e-mail addresses (shared-string values) ---
$ws.Range("D2").Value = "murbina@yopmail.com"
$ws.Range("D3").Value = "eurbina@yopmail.com"
$ws.Range("D4").Value = "jmagallanes@yopmail.com"

# --- Widen column D ---
$ws.Columns.Item(4).ColumnWidth = 25.15

# --- Move / persist the user's current selection ---
$ws.Range("D9").Select()
